# "updated params and drawer locking"
#
# Applies the changes described by the diff to Sheet1 of the workbook:
#  - B22 formula updated to include B57 and B56/2 terms
#  - D57 category changed from "y" to "new" (drawer-locking related row)
#  - B60 value bumped from 45 to 50
#  - B81 formula switched from a CEILING() of a few params to a FLOOR()
#    based on B78/B80/B79/B60
#  - B82 changed from a hard-coded constant to a formula (B57+B62)
#  - the view is scrolled/selected near the bottom of the sheet (A60 / B61)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- parameter / formula updates -------------------------------------------

# B22: B77+B23+B1  ->  B77+B23+B57-B56/2-B1
$ws.Range("B22").Formula = "=B77+B23+B57-B56/2-B1"

# D57: "y" -> "new"  (drawer locking category)
$ws.Range("D57").Value = "new"

# B60: 45 -> 50
$ws.Range("B60").Value2 = 50

# B81: CEILING(B58+B60+25,1) -> FLOOR(B78-B80-B79-B60,1)
# (must be set after B60 is updated so the cached value matches)
$ws.Range("B81").Formula = "=FLOOR(B78-B80-B79-B60,1)"

# B82: constant 50 -> formula B57+B62
$ws.Range("B82").Formula = "=B57+B62"

# --- view state --------------------------------------------------------
# topLeftCell A42 -> A60, selection D53 -> B61
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B61").Select()
